$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Girokonto")
$ws2 = $wb.Worksheets.Item("Kreditkarte")

# Update the long hint text in D4 on Girokonto: Strg+Enter -> Alt+Enter
$ws1.Range("D4").Value = "Lange Kategorienamen können per erzwungenen Zeilenumbruch (Alt+Enter) in den Diagramm umgebrochen dargestellt werden."

# Kreditkarte keeps its old selection (B12) but is no longer the active/front tab
$ws2.Range("B12").Select()

# Girokonto becomes the active sheet/tab, with D4 selected (and no frozen/scrolled topLeftCell)
$ws1.Activate()
$ws1.Range("D4").Select()
